# Generate Report for Handoff
# A new handoff cycle was generated for 049a2d9b-...md (status -> "Ready for
# handoff", new handoff datetime), while a6d6703e-...md keeps its prior
# "Handed back: in sync with en-US" status. The two source rows also swap
# display order/position across all three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-36-20 10:36:18"

$ov.Range("A3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-37-20 10:37:16"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-20 10:36:15"
$zh.Range("F2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.md"
$zh.Range("G2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-20 10:36:40"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 10:37:12"
$zh.Range("F3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md"
$zh.Range("G3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-20 10:36:40"
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/54d7f3c61687d4a869021ce2c659e3f5be5d4487/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fa416c202a5eef3f4e1f16f937b29386a4b87f31/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/86b61339b5880fbd1d3b2b675ed94527ffd679f9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/54d7f3c61687d4a869021ce2c659e3f5be5d4487/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fa416c202a5eef3f4e1f16f937b29386a4b87f31/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/86b61339b5880fbd1d3b2b675ed94527ffd679f9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf"
$de.Range("E2").Value = "2016-03-20 10:36:18"
$de.Range("F2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.md"
$de.Range("G2").Value = "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf"
$de.Range("H2").Value = "2016-03-20 10:36:46"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 10:37:16"
$de.Range("F3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md"
$de.Range("G3").Value = "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf"
$de.Range("H3").Value = "2016-03-20 10:36:46"
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9665284558c08d724da120f51941c0aacbcd802c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5a5ffdf6df8065e076132742f658003b1e98157d/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/a6d6703e-4025-4230-ad00-6c893f049051.md", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7784400865c8bc8e3ad323d41bfbcb1e94db31d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf", [Type]::Missing, [Type]::Missing, "a6d6703e-4025-4230-ad00-6c893f049051.c60d92997ecc0806b6d38610f574fcddea34b94b.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9665284558c08d724da120f51941c0aacbcd802c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5a5ffdf6df8065e076132742f658003b1e98157d/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3375b8c987af8991c5baad40dec9323fbea7b0e5/e2e/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7784400865c8bc8e3ad323d41bfbcb1e94db31d4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf", [Type]::Missing, [Type]::Missing, "049a2d9b-6e27-4c09-9e57-ca0ded622ef7.f88a9ebb3d84892f83454be91a3df842bdad22d4.de-de.xlf") | Out-Null
